$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at position 10 (new weekly price report) ---
$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value = "La Araucanía"
$ws.Cells.Item(10, 4).Value = (Get-Date -Year 2023 -Month 4 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100108
$ws.Cells.Item(10, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(10, 9).Value = 100108003
$ws.Cells.Item(10, 10).Value = "Maracuyá"
$ws.Cells.Item(10, 11).Value = "Sin especificar"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 40
$ws.Cells.Item(10, 14).Value = 45000
$ws.Cells.Item(10, 15).Value = 45000
$ws.Cells.Item(10, 16).Value = 45000
$ws.Cells.Item(10, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(10, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(10, 19).Value = 2500
$ws.Cells.Item(10, 20).Value = 18

# --- Insert a second new row at position 82 (after the first insert shifted
#     everything down by one) for another new weekly price report ---
$ws.Rows.Item(82).Insert()

$ws.Cells.Item(82, 1).Value = 10
$ws.Cells.Item(82, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(82, 3).Value = "La Araucanía"
$ws.Cells.Item(82, 4).Value = (Get-Date -Year 2023 -Month 4 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(82, 5).Value = 9
$ws.Cells.Item(82, 6).Value = "Fruta"
$ws.Cells.Item(82, 7).Value = 100108
$ws.Cells.Item(82, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(82, 9).Value = 100108003
$ws.Cells.Item(82, 10).Value = "Maracuyá"
$ws.Cells.Item(82, 11).Value = "Sin especificar"
$ws.Cells.Item(82, 12).Value = "Primera"
$ws.Cells.Item(82, 13).Value = 50
$ws.Cells.Item(82, 14).Value = 45000
$ws.Cells.Item(82, 15).Value = 45000
$ws.Cells.Item(82, 16).Value = 45000
$ws.Cells.Item(82, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(82, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(82, 19).Value = 2500
$ws.Cells.Item(82, 20).Value = 18
